# Update "想去人数" (want-to-go count) figures across sheets to match
# the freshly generated gh-pages output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1100
$ws1.Range("F8").Value = 391
$ws1.Range("F9").Value = 1010
$ws1.Range("F13").Value = 154
$ws1.Range("F14").Value = 12614
$ws1.Range("F15").Value = 5202
$ws1.Range("F16").Value = 5517

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 30

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1100
$ws4.Range("F9").Value = 391
$ws4.Range("F10").Value = 1010
$ws4.Range("F14").Value = 154
$ws4.Range("F15").Value = 12614
$ws4.Range("F16").Value = 30
$ws4.Range("F18").Value = 5202
$ws4.Range("F19").Value = 5517
